$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 1372.875
$ws.Range("I121").Value = 995
$ws.Range("J121").Value = 1599.6
$ws.Range("K121").Value = 2985
$ws.Range("L121").Value = 4798.799999999999
$ws.Range("M121").Value = -1238
$ws.Range("N121").Value = -8292.799999999999
$ws.Range("H130").Value = 250044240
$ws.Range("J130").Value = 250044240
$ws.Range("L130").Value = 250044240
$ws.Range("N130").Value = -250054280
$ws.Range("H132").Value = 6579.622
$ws.Range("I132").Value = 5863.6787
$ws.Range("K132").Value = 17591.0361
$ws.Range("M132").Value = -15061.0361
$ws.Range("H137").Value = 2003.6774
$ws.Range("I137").Value = 2130.56
$ws.Range("J137").Value = 1475
$ws.Range("K137").Value = 6391.68
$ws.Range("L137").Value = 4425
$ws.Range("M137").Value = -3841.68
$ws.Range("N137").Value = -9525

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3076.9412
$ws.Range("I2").Value = 2053.8667
$ws.Range("K2").Value = 2053.8667
$ws.Range("M2").Value = -1940.8667
$ws.Range("H116").Value = 3076.9412
$ws.Range("I116").Value = 2053.8667
$ws.Range("K116").Value = 2053.8667
$ws.Range("M116").Value = 240.1333
$ws.Range("H132").Value = 3426.2542
$ws.Range("I132").Value = 1336.5952
$ws.Range("J132").Value = 8588.941000000001
$ws.Range("K132").Value = 4009.7856
$ws.Range("L132").Value = 25766.823
$ws.Range("M132").Value = -1479.7856
$ws.Range("N132").Value = -30826.823

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3076.9412
$ws.Range("I3").Value = 2053.8667
$ws.Range("K3").Value = 2053.8667
$ws.Range("M3").Value = -1939.8667
$ws.Range("H15").Value = 5000
$ws.Range("J15").Value = 5000
$ws.Range("L15").Value = 5000
$ws.Range("N15").Value = -5454
$ws.Range("H134").Value = 5288.787
$ws.Range("I134").Value = 2455.476
$ws.Range("J134").Value = 7577.231
$ws.Range("K134").Value = 7366.428
$ws.Range("L134").Value = 22731.693
$ws.Range("M134").Value = -4831.428
$ws.Range("N134").Value = -27801.693

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6804577
$ws.Range("I31").Value = 1344.1
$ws.Range("J31").Value = 37041170
$ws.Range("K31").Value = 1344.1
$ws.Range("L31").Value = 37041170
$ws.Range("M31").Value = -1049.1
$ws.Range("N31").Value = -37041760
$ws.Range("H34").Value = 6804577
$ws.Range("I34").Value = 1344.1
$ws.Range("J34").Value = 37041170
$ws.Range("K34").Value = 1344.1
$ws.Range("L34").Value = 37041170
$ws.Range("M34").Value = -1142.1
$ws.Range("N34").Value = -37041574
$ws.Range("H62").Value = 19793.924
$ws.Range("J62").Value = 23626.666
$ws.Range("L62").Value = 23626.666
$ws.Range("N62").Value = -24874.666
$ws.Range("H65").Value = 19793.924
$ws.Range("J65").Value = 23626.666
$ws.Range("L65").Value = 118133.33
$ws.Range("N65").Value = -124373.33
$ws.Range("H99").Value = 2678.5
$ws.Range("I99").Value = 4400
$ws.Range("J99").Value = 2432.5715
$ws.Range("K99").Value = 4400
$ws.Range("L99").Value = 2432.5715
$ws.Range("M99").Value = -2902
$ws.Range("N99").Value = -5428.5715
$ws.Range("H126").Value = 2678.5
$ws.Range("I126").Value = 4400
$ws.Range("J126").Value = 2432.5715
$ws.Range("K126").Value = 13200
$ws.Range("L126").Value = 7297.7145
$ws.Range("M126").Value = -10730
$ws.Range("N126").Value = -12237.7145
$ws.Range("H132").Value = 2612.9143
$ws.Range("I132").Value = 1535.9131
$ws.Range("J132").Value = 4677.1665
$ws.Range("K132").Value = 4607.7393
$ws.Range("L132").Value = 14031.4995
$ws.Range("M132").Value = -2077.7393
$ws.Range("N132").Value = -19091.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2759.2
$ws.Range("I80").Value = 975.25
$ws.Range("J80").Value = 3205.1875
$ws.Range("K80").Value = 2925.75
$ws.Range("L80").Value = 9615.5625
$ws.Range("M80").Value = -1989.75
$ws.Range("N80").Value = -11487.5625
$ws.Range("H83").Value = 2759.2
$ws.Range("I83").Value = 975.25
$ws.Range("J83").Value = 3205.1875
$ws.Range("K83").Value = 8777.25
$ws.Range("L83").Value = 28846.6875
$ws.Range("M83").Value = -4097.25
$ws.Range("N83").Value = -38206.6875
$ws.Range("H113").Value = 632.9268
$ws.Range("I113").Value = 492
$ws.Range("K113").Value = 1476
$ws.Range("M113").Value = 694
$ws.Range("H122").Value = 2563.2
$ws.Range("J122").Value = 3237.5186
$ws.Range("L122").Value = 29137.6674
$ws.Range("N122").Value = -34037.6674
$ws.Range("H131").Value = 323441.78
$ws.Range("I131").Value = 769741.25
$ws.Range("J131").Value = 1114.3889
$ws.Range("K131").Value = 2309223.75
$ws.Range("L131").Value = 3343.1667
$ws.Range("M131").Value = -2304183.75
$ws.Range("N131").Value = -13423.1667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H127").Value = 25000
$ws.Range("J127").Value = 25000
$ws.Range("L127").Value = 25000
$ws.Range("N127").Value = -34920
$ws.Range("H132").Value = 2427.4324
$ws.Range("I132").Value = 1763.1428
$ws.Range("J132").Value = 3299.3125
$ws.Range("K132").Value = 5289.428400000001
$ws.Range("L132").Value = 9897.9375
$ws.Range("M132").Value = -2759.428400000001
$ws.Range("N132").Value = -14957.9375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 1597.5
$ws.Range("I5").Value = 796.6667
$ws.Range("K5").Value = 796.6667
$ws.Range("M5").Value = -683.6667
$ws.Range("H7").Value = 3750
$ws.Range("I7").Value = 4000
$ws.Range("J7").Value = 3500
$ws.Range("K7").Value = 4000
$ws.Range("L7").Value = 3500
$ws.Range("M7").Value = -3888
$ws.Range("N7").Value = -3724
$ws.Range("H40").Value = 50003788
$ws.Range("I40").Value = 90911630
$ws.Range("K40").Value = 90911630
$ws.Range("M40").Value = -90911494
$ws.Range("H46").Value = 556827.4
$ws.Range("I46").Value = 850.25
$ws.Range("J46").Value = 715678
$ws.Range("K46").Value = 850.25
$ws.Range("L46").Value = 715678
$ws.Range("M46").Value = -662.25
$ws.Range("N46").Value = -716054
$ws.Range("H126").Value = 3750
$ws.Range("I126").Value = 4000
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 12000
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -9530
$ws.Range("N126").Value = -15440
$ws.Range("H127").Value = 29636.666
$ws.Range("J127").Value = 29636.666
$ws.Range("L127").Value = 29636.666
$ws.Range("N127").Value = -39556.666
$ws.Range("H130").Value = 28590
$ws.Range("J130").Value = 28590
$ws.Range("L130").Value = 28590
$ws.Range("N130").Value = -38630

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 600
$ws.Range("I26").Value = 600
$ws.Range("K26").Value = 600
$ws.Range("M26").Value = -307
$ws.Range("H41").Value = 26594.25
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 26594.25
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 26594.25
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -27374.25
$ws.Range("H122").Value = 5766.8887
$ws.Range("I122").Value = 5003.769
$ws.Range("J122").Value = 7751
$ws.Range("K122").Value = 15011.307
$ws.Range("L122").Value = 23253
$ws.Range("M122").Value = -12561.307
$ws.Range("N122").Value = -28153
$ws.Range("H126").Value = 3059.9473
$ws.Range("I126").Value = 2819.5
$ws.Range("J126").Value = 3472.1428
$ws.Range("K126").Value = 8458.5
$ws.Range("L126").Value = 10416.4284
$ws.Range("M126").Value = -5988.5
$ws.Range("N126").Value = -15356.4284
$ws.Range("H132").Value = 2761.7307
$ws.Range("I132").Value = 2133.6667
$ws.Range("J132").Value = 5399.6
$ws.Range("K132").Value = 6401.000100000001
$ws.Range("L132").Value = 16198.8
$ws.Range("M132").Value = -3871.000100000001
$ws.Range("N132").Value = -21258.8
